$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row (based on column A, which is always populated).
$lastRow = $ws.Cells($ws.Rows.Count, 1).End(-4162).Row

# Update the "Förändrad" (changed) date in column C for every existing data row
# (rows 2..lastRow) from 45202 to 45203.
$ws.Range("C2:C$lastRow").Value2 = 45203

# Touch the row height of the last existing row so it carries an explicit
# custom height, matching what Excel does when a row is appended after it.
$ws.Rows.Item($lastRow).RowHeight = 15

# Append the new record as the next row.
$newRow = $lastRow + 1

$ws.Cells.Item($newRow, 1).Value2 = "A 47347-2023"

$ws.Cells.Item($newRow, 2).Value2 = 45202
$ws.Cells.Item($newRow, 2).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item($newRow, 3).Value2 = 45203
$ws.Cells.Item($newRow, 3).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item($newRow, 4).Value2 = "DALARNAS LÄN"
$ws.Cells.Item($newRow, 5).Value2 = "MORA"
$ws.Cells.Item($newRow, 6).Value2 = "Bergvik skog öst AB"
$ws.Cells.Item($newRow, 7).Value2 = 5.1
$ws.Cells.Item($newRow, 8).Value2 = 0
$ws.Cells.Item($newRow, 9).Value2 = 0
$ws.Cells.Item($newRow, 10).Value2 = 0
$ws.Cells.Item($newRow, 11).Value2 = 0
$ws.Cells.Item($newRow, 12).Value2 = 0
$ws.Cells.Item($newRow, 13).Value2 = 0
$ws.Cells.Item($newRow, 14).Value2 = 0
$ws.Cells.Item($newRow, 15).Value2 = 0
$ws.Cells.Item($newRow, 16).Value2 = 0
$ws.Cells.Item($newRow, 17).Value2 = 0

# Column R keeps the wrap-text style seen on every other row, left blank.
$ws.Cells.Item($newRow, 18).WrapText = $true
